$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered) from A1 onto F1
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F1").Value = "scenario"

for ($i = 2; $i -le 101; $i++) {
    $ws.Cells.Item($i, 6).Value = "S6"
}
